# Bump the "Förändrad" (changed) date in column C by one day for every data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C header is in row 1 ("Förändrad"); data rows start at row 2.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 408 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46061) {
        $cell.Value2 = 46062
    }
}
